$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new row for the latest week at row 3,
# pushing all existing data rows (previously 3-10) down to 4-11.
$ws.Rows("3:3").Insert()

# Populate the new week's row with its data.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44453
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 850
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Región de Arica y Parinacota"
$ws.Range("P3").Value = 850
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
